$wb = $excel.ActiveWorkbook

# --- Sheet: P_valores ---
$wsP = $wb.Worksheets.Item("P_valores")

$wsP.Range("C2").Value = 0.9178988381992794
$wsP.Range("D2").Value = 0.9960422874864701
$wsP.Range("E2").Value = 0.7310413714076835
$wsP.Range("F2").Value = 0.5318041844337564

$wsP.Range("B3").Value = 0.9178988381992794
$wsP.Range("D3").Value = 0.9003135897864745
$wsP.Range("E3").Value = 0.5583582407298047
$wsP.Range("F3").Value = 0.3838788986999662

$wsP.Range("B4").Value = 0.9960422874864701
$wsP.Range("C4").Value = 0.9003135897864745
$wsP.Range("E4").Value = 0.5920910891106557
$wsP.Range("F4").Value = 0.3019004315381617

$wsP.Range("B5").Value = 0.7310413714076835
$wsP.Range("C5").Value = 0.5583582407298047
$wsP.Range("D5").Value = 0.5920910891106557
$wsP.Range("F5").Value = 0.1411561762585349

$wsP.Range("B6").Value = 0.5318041844337564
$wsP.Range("C6").Value = 0.3838788986999662
$wsP.Range("D6").Value = 0.3019004315381617
$wsP.Range("E6").Value = 0.1411561762585349

# --- Sheet: Estadisticos_DM ---
$wsE = $wb.Worksheets.Item("Estadisticos_DM")

$wsE.Range("C2").Value = -0.1049569952345589
$wsE.Range("D2").Value = -0.005049575427194796
$wsE.Range("E2").Value = 0.3506904956766931
$wsE.Range("F2").Value = -0.6411141255092242

$wsE.Range("B3").Value = 0.1049569952345589
$wsE.Range("D3").Value = 0.1275570710381784
$wsE.Range("E3").Value = 0.5995915562271927
$wsE.Range("F3").Value = -0.898952806297762

$wsE.Range("B4").Value = 0.005049575427194796
$wsE.Range("C4").Value = -0.1275570710381784
$wsE.Range("E4").Value = 0.5483466213747694
$wsE.Range("F4").Value = -1.071889061149453

$wsE.Range("B5").Value = -0.3506904956766931
$wsE.Range("C5").Value = -0.5995915562271927
$wsE.Range("D5").Value = -0.5483466213747694
$wsE.Range("F5").Value = -1.559652995885853

$wsE.Range("B6").Value = 0.6411141255092242
$wsE.Range("C6").Value = 0.898952806297762
$wsE.Range("D6").Value = 1.071889061149453
$wsE.Range("E6").Value = 1.559652995885853
